# Update automàtic: dades i banners [2026-02-12 23:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-valued cells: force Text format first so Excel does not
# auto-convert the literal "NN%" string into a numeric percentage value,
# then restore General formatting once the literal text is stored.
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "40%"
$ws.Range("H6").NumberFormat = "General"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "49%"
$ws.Range("H11").NumberFormat = "General"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "70%"
$ws.Range("H12").NumberFormat = "General"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "51%"
$ws.Range("H13").NumberFormat = "General"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "51%"
$ws.Range("H15").NumberFormat = "General"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "74%"
$ws.Range("H17").NumberFormat = "General"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "86%"
$ws.Range("H20").NumberFormat = "General"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "62%"
$ws.Range("H29").NumberFormat = "General"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "46%"
$ws.Range("H31").NumberFormat = "General"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "53%"
$ws.Range("H33").NumberFormat = "General"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "57%"
$ws.Range("H40").NumberFormat = "General"

# Remaining cell updates (timestamps, pressures, temperatures, etc.)
$ws.Range("E2").Value = "2026-02-12 23:18:31"
$ws.Range("E3").Value = "2026-02-12 23:18:33"
$ws.Range("E4").Value = "2026-02-12 23:18:36"
$ws.Range("J4").Value = "1000.0 hPa"
$ws.Range("O4").Value = "16.0 °C"
$ws.Range("E5").Value = "2026-02-12 23:18:39"
$ws.Range("E6").Value = "2026-02-12 23:18:41"
$ws.Range("J6").Value = "999.9 hPa"
$ws.Range("O6").Value = "15.6 °C"
$ws.Range("E7").Value = "2026-02-12 23:18:44"
$ws.Range("J7").Value = "1002.4 hPa"
$ws.Range("E8").Value = "2026-02-12 23:18:46"
$ws.Range("E9").Value = "2026-02-12 23:18:49"
$ws.Range("K9").Value = "13.3 MJ/m2"
$ws.Range("E10").Value = "2026-02-12 23:18:51"
$ws.Range("N10").Value = "7.6 °C 22:54 TU"
$ws.Range("O10").Value = "14.5 °C"
$ws.Range("E11").Value = "2026-02-12 23:18:54"
$ws.Range("O11").Value = "8.9 °C"
$ws.Range("E12").Value = "2026-02-12 23:18:56"
$ws.Range("E13").Value = "2026-02-12 23:18:59"
$ws.Range("J13").Value = "1002.6 hPa"
$ws.Range("O13").Value = "7.4 °C"
$ws.Range("E14").Value = "2026-02-12 23:19:01"
$ws.Range("O14").Value = "16.8 °C"
$ws.Range("E15").Value = "2026-02-12 23:19:04"
$ws.Range("E16").Value = "2026-02-12 23:19:06"
$ws.Range("E17").Value = "2026-02-12 23:19:09"
$ws.Range("E18").Value = "2026-02-12 23:19:11"
$ws.Range("J18").Value = "1000.3 hPa"
$ws.Range("N18").Value = "8.0 °C 22:54 TU"
$ws.Range("O18").Value = "16.3 °C"
$ws.Range("E19").Value = "2026-02-12 23:19:14"
$ws.Range("O19").Value = "7.8 °C"
$ws.Range("E20").Value = "2026-02-12 23:19:17"
$ws.Range("E21").Value = "2026-02-12 23:19:19"
$ws.Range("J21").Value = "1003.1 hPa"
$ws.Range("O21").Value = "8.9 °C"
$ws.Range("E22").Value = "2026-02-12 23:19:22"
$ws.Range("E23").Value = "2026-02-12 23:19:24"
$ws.Range("E24").Value = "2026-02-12 23:19:26"
$ws.Range("E25").Value = "2026-02-12 23:19:29"
$ws.Range("O25").Value = "-2.0 °C"
$ws.Range("E26").Value = "2026-02-12 23:19:31"
$ws.Range("J26").Value = "999.7 hPa"
$ws.Range("E27").Value = "2026-02-12 23:19:34"
$ws.Range("E28").Value = "2026-02-12 23:19:37"
$ws.Range("J28").Value = "999.8 hPa"
$ws.Range("O28").Value = "13.5 °C"
$ws.Range("E29").Value = "2026-02-12 23:19:39"
$ws.Range("N29").Value = "5.4 °C 22:55 TU"
$ws.Range("O29").Value = "13.7 °C"
$ws.Range("E30").Value = "2026-02-12 23:19:42"
$ws.Range("J30").Value = "1000.1 hPa"
$ws.Range("N30").Value = "6.5 °C 22:32 TU"
$ws.Range("O30").Value = "11.7 °C"
$ws.Range("E31").Value = "2026-02-12 23:19:44"
$ws.Range("J31").Value = "999.5 hPa"
$ws.Range("E32").Value = "2026-02-12 23:19:47"
$ws.Range("E33").Value = "2026-02-12 23:19:49"
$ws.Range("J33").Value = "1002.2 hPa"
$ws.Range("N33").Value = "2.4 °C 22:56 TU"
$ws.Range("E34").Value = "2026-02-12 23:19:52"
$ws.Range("E35").Value = "2026-02-12 23:19:55"
$ws.Range("E36").Value = "2026-02-12 23:19:57"
$ws.Range("J36").Value = "1000.4 hPa"
$ws.Range("E37").Value = "2026-02-12 23:20:00"
$ws.Range("J37").Value = "1001.3 hPa"
$ws.Range("N37").Value = "2.4 °C 22:34 TU"
$ws.Range("O37").Value = "9.4 °C"
$ws.Range("E38").Value = "2026-02-12 23:20:02"
$ws.Range("N38").Value = "12.0 °C 22:54 TU"
$ws.Range("E39").Value = "2026-02-12 23:20:05"
$ws.Range("E40").Value = "2026-02-12 23:20:08"
$ws.Range("J40").Value = "1003.9 hPa"
$ws.Range("N40").Value = "3.1 °C 22:55 TU"
$ws.Range("O40").Value = "9.0 °C"
$ws.Range("E41").Value = "2026-02-12 23:20:10"
$ws.Range("J41").Value = "1005.9 hPa"
$ws.Range("E42").Value = "2026-02-12 23:20:13"
$ws.Range("O42").Value = "13.5 °C"
$ws.Range("E43").Value = "2026-02-12 23:20:15"
$ws.Range("O43").Value = "11.8 °C"
$ws.Range("E44").Value = "2026-02-12 23:20:17"
$ws.Range("N44").Value = "-6.0 °C 22:55 TU"
$ws.Range("E45").Value = "2026-02-12 23:20:20"
$ws.Range("J45").Value = "1005.6 hPa"
$ws.Range("N45").Value = "1.4 °C 22:36 TU"
$ws.Range("O45").Value = "6.6 °C"
$ws.Range("E46").Value = "2026-02-12 23:20:23"
$ws.Range("O46").Value = "15.6 °C"
